$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.004, -0.672, -0.433, -65, -18, -100),
    @(0.171, 0.13, -0.307, -115, 44, -7),
    @(-0.323, 0.167, -0.212, -22, -46, -82),
    @(0.475, -0.104, -0.651, 100, 62, -112),
    @(0.217, 0.255, 1.019, -57, 39, 89),
    @(-0.199, 0.08, 0.4, -164, 68, 13),
    @(-0.122, 1.004, 0.032, -91, 3, -78),
    @(-0.484, 0.619, 0.894, -126, 63, 1),
    @(0.628, 0.285, 0.184, -88, -38, 6),
    @(0.07199999999999999, -0.654, -0.33, 84, 11, 168)
)

for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$r-1][$c-1]
    }
}
